$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.888.50'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.808.85'
$ws.Range('E3').Value = '  -3.36%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.60'
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('E6').Value = '  -1.71%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '39.37'
$ws.Range('E8').Value = '  -7.82%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.317'
$ws.Range('E9').Value = '  +1.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0994'
$ws.Range('E11').Value = '  -1.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.070.79'
$ws.Range('E12').Value = '  -3.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.823.14'
$ws.Range('E13').Value = '  -2.59%  '
$ws.Range('E14').Value = '  -4.26%  '
$ws.Range('E15').Value = '  -7.17%  '
$ws.Range('E16').Value = '  -5.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '34.835.71'
$ws.Range('E17').Value = '  -2.74%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.02'
$ws.Range('E18').Value = '  -2.60%  '
$ws.Range('E19').Value = '  -3.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '238.05'
$ws.Range('E20').Value = '  -4.65%  '
$ws.Range('E21').Value = '  -6.00%  '
$ws.Range('E22').Value = '  -4.36%  '
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('E24').Value = '  -1.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '172.64'
$ws.Range('E25').Value = '  +0.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.79'
$ws.Range('E26').Value = '  -4.05%  '
$ws.Range('E27').Value = '  -4.72%  '
$ws.Range('E28').Value = '  -3.58%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.55'
$ws.Range('E29').Value = '  +7.24%  '
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('E31').Value = '  -0.74%  '
$ws.Range('E32').Value = '  -1.95%  '
$ws.Range('E33').Value = '  -4.03%  '
$ws.Range('E34').Value = '  -8.81%  '
$ws.Range('E35').Value = '  +4.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.683'
$ws.Range('E36').Value = '  -1.38%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '91.08'
$ws.Range('E37').Value = '  -8.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.31'
$ws.Range('E38').Value = '  +3.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.306.51'
$ws.Range('E39').Value = '  -4.66%  '
$ws.Range('E40').Value = '  -3.73%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.45'
$ws.Range('E41').Value = '  -1.40%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.958'
$ws.Range('E42').Value = '  -7.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '14.23'
$ws.Range('E43').Value = '  -6.00%  '
$ws.Range('E44').Value = '  -13.11%  '
$ws.Range('E45').Value = '  -5.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.12'
$ws.Range('E46').Value = '  -3.21%  '
$ws.Range('E47').Value = '  -2.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.991.83'
$ws.Range('E48').Value = '  -2.45%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.01'
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0670'
$ws.Range('E50').Value = '  +6.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '98.63'
$ws.Range('E51').Value = '  -6.50%  '
